$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats/styles from the (now-shifted) F:M data range into the
# newly inserted D:E columns so the new cells match the rest of the row
$ws.Range("F5:M102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 14710000
$ws.Range("E8").Value = 15366000
$ws.Range("D9").Value = 12098000
$ws.Range("E9").Value = 12669000
$ws.Range("D10").Value = 2612000
$ws.Range("E10").Value = 2697000
$ws.Range("D12").Value = 344000
$ws.Range("E12").Value = 354000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 65000
$ws.Range("E14").Value = 73000
$ws.Range("D15").Value = 29000
$ws.Range("E15").Value = 20000
$ws.Range("D17").Value = 13784000
$ws.Range("E17").Value = 14319000
$ws.Range("D18").Value = 926000
$ws.Range("E18").Value = 1047000
$ws.Range("D20").Value = 38000
$ws.Range("E20").Value = 31000
$ws.Range("D21").Value = 1132000
$ws.Range("E21").Value = 1218000
$ws.Range("D22").Value = 64000
$ws.Range("E22").Value = 71000
$ws.Range("D23").Value = 900000
$ws.Range("E23").Value = 1007000
$ws.Range("D24").Value = 118000
$ws.Range("E24").Value = -544000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 782000
$ws.Range("E26").Value = 1551000
$ws.Range("D27").Value = 782000
$ws.Range("E27").Value = 1551000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 21000
$ws.Range("E29").Value = -100000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -38000
$ws.Range("E32").Value = -31000
$ws.Range("D33").Value = 803000
$ws.Range("E33").Value = 1451000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 803000
$ws.Range("E35").Value = 1451000
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 3367000
$ws.Range("E41").Value = 5166000
$ws.Range("D42").Value = 425000
$ws.Range("E42").Value = 711000
$ws.Range("D43").Value = 8211000
$ws.Range("E43").Value = 8003000
$ws.Range("D44").Value = 5649000
$ws.Range("E44").Value = 6062000
$ws.Range("D45").Value = 1284000
$ws.Range("E45").Value = 1445000
$ws.Range("D46").Value = 18936000
$ws.Range("E46").Value = 21387000
$ws.Range("D47").Value = 915000
$ws.Range("E47").Value = 989000
$ws.Range("D48").Value = 2312000
$ws.Range("E48").Value = 2198000
$ws.Range("D49").Value = 7067000
$ws.Range("E49").Value = 5968000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3260000
$ws.Range("E52").Value = 4080000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 32490000
$ws.Range("E54").Value = 34622000
$ws.Range("D57").Value = 14572000
$ws.Range("E57").Value = 14816000
$ws.Range("D58").Value = 297000
$ws.Range("E58").Value = 1463000
$ws.Range("D59").Value = 9330000
$ws.Range("E59").Value = 8852000
$ws.Range("D60").Value = 24199000
$ws.Range("E60").Value = 25131000
$ws.Range("D61").Value = 4706000
$ws.Range("E61").Value = 4524000
$ws.Range("D62").Value = 5422000
$ws.Range("E62").Value = 5606000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 34327000
$ws.Range("E66").Value = 35261000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1431000
$ws.Range("E72").Value = -473000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -1837000
$ws.Range("E76").Value = -639000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = 803000
$ws.Range("E81").Value = 1451000
$ws.Range("D83").Value = 168000
$ws.Range("E83").Value = 140000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 862000
$ws.Range("E89").Value = 968000
$ws.Range("D91").Value = -189000
$ws.Range("E91").Value = -187000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -318000
$ws.Range("E94").Value = 87000
$ws.Range("D96").Value = -249000
$ws.Range("E96").Value = -219000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2343000
$ws.Range("E100").Value = -2084000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -1799000
$ws.Range("E102").Value = -1029000
